# Weekly fruit/vegetable price update: insert a new week's record as the
# new row 97 (pushing the existing rows 97-134 down to 98-135), matching
# the rest of the sheet's pattern of one row per market report.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 97; this shifts rows
# 97-134 down to 98-135 and extends the sheet dimension automatically.
$ws.Rows(97).Insert()

# Populate the newly inserted row 97 with the new weekly record. All
# fields mirror the record that used to be at row 97 (now row 98) except
# for the reporting date (D) and the traded volume (J).
$ws.Range("A97").Value = 10
$ws.Range("B97").Value = "Vega Modelo de Temuco"
$ws.Range("C97").Value = "La Araucanía"
$ws.Range("D97").Value = 45146
$ws.Range("E97").Value = 9
$ws.Range("F97").Value = 100112010
$ws.Range("G97").Value = "Achicoria"
$ws.Range("H97").Value = "Sin especificar"
$ws.Range("I97").Value = "Primera"
$ws.Range("J97").Value = 65
$ws.Range("K97").Value = 10000
$ws.Range("L97").Value = 10000
$ws.Range("M97").Value = 10000
$ws.Range("N97").Value = "$/caja 18 unidades"
$ws.Range("O97").Value = "Región Metropolitana"
$ws.Range("P97").Value = 556
$ws.Range("Q97").Value = 18
$ws.Range("R97").Value = "Hortaliza"
